# Implement checkout functionality: add CheckoutConfirmationPage, enhance shopping
# cart tests, and refactor UI actions for better reporting.
#
# Concretely (as reflected in the AddProduct / DeleteProduct test-data sheets):
#   - Disable ("No") the ExecutionFlag for every AddProduct test case except
#     TC0001, which stays enabled ("Yes").
#   - The AddProduct sheet becomes the active/selected tab (with cell D13
#     selected), while DeleteProduct is no longer the active tab (its own
#     selection of D12 is left untouched).

$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("AddProduct")
$wsDelete = $wb.Worksheets.Item("DeleteProduct")

# Turn off ExecutionFlag ("Yes" -> "No") for TC0002..TC0007 (rows 3 through 8).
# Row 2 (TC0001) is left as "Yes".
foreach ($r in 3..8) {
    $wsAdd.Cells.Item($r, 1).Value = "No"
}

# Make AddProduct the active sheet/tab with D13 selected.
[void]$wsAdd.Activate()
[void]$wsAdd.Range("D13").Select()

# DeleteProduct keeps its own prior selection (D12), it's just no longer the
# active tab once AddProduct has been activated above.
[void]$wsDelete.Range("D12").Select()
[void]$wsAdd.Activate()
